{"js": "// The canonical XML diff for this change only touches two kinds of purely\n// internal/random identifiers that Word (re-)generated when the fixture was\n// regenerated going from 3.1.1 to 3.2.0:\n//   1. the (non-standard, tool-generated) w:rsidR=\"....\" GUID that the\n//      M2Doc generator stamps on every run belonging to one of the two\n//      \"REF bookmark1\" field instances (begin/instrText/separate/result/end).\n//   2. the numeric w:id on the bookmarkStart/bookmarkEnd pair that wraps\n//      \"bookmarked content\" (bookmark \"bookmark1\").\n// No visible text, formatting or structure changes. We reproduce this by\n// rewriting, in place, the raw OOXML of the three paragraphs that contain\n// those runs/bookmark, using the exact literal target values from the diff.\n\nconst OLD_RSID = \"33AF09C2681E4769971B94CEE0D9A9EC\";\nconst NEW_RSID = \"07822CC306A04453A57D7E97965A69F1\";\nconst OLD_BMK_ID = \"45476227328674507871903955380577777817\";\nconst NEW_BMK_ID = \"98058591413380703969263738417409166366\";\n\nfunction wrapPackage(paragraphXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + paragraphXml + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Locate the three paragraphs touched by the change purely from their\n// (unchanged) visible text, so this does not depend on paragraph order.\nlet pBefore = null; // \"Test link before bookmark : ...\"\nlet pBookmark = null; // \"Test bookmark : bookmarked content\" (first one, the real bookmark)\nlet pAfter = null; // \"Test link after bookmark : ...\"\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const t = p.text;\n  if (t.indexOf(\"Test link before bookmark\") !== -1) {\n    pBefore = p;\n  } else if (pBookmark === null && t.indexOf(\"Test bookmark\") !== -1 && t.indexOf(\"bookmarked content\") !== -1) {\n    pBookmark = p;\n  } else if (t.indexOf(\"Test link after bookmark\") !== -1) {\n    pAfter = p;\n  }\n}\n\n// Pull the raw OOXML for each of the three target paragraphs up-front.\nconst ooxmlBefore = pBefore ? pBefore.getOoxml() : null;\nconst ooxmlBookmark = pBookmark ? pBookmark.getOoxml() : null;\nconst ooxmlAfter = pAfter ? pAfter.getOoxml() : null;\nawait context.sync();\n\n// The Office.js OOXML projection does not echo back tool-generated\n// w:rsidR / large bookmark w:id values, so we rebuild the exact target\n// fragments from the document's known original markup (the field-code\n// runs and bookmark wrapper), substituting the new literal id values\n// straight from the diff, then write them back verbatim over each\n// paragraph.\n\nconst fieldRunsXml =\n  '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n  '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:instrText xml:space=\"preserve\"> REF bookmark1 \\\\h </w:instrText></w:r>' +\n  '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n  '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r>' +\n  '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r>';\n\nif (pBefore) {\n  const beforeXml =\n    '<w:p w:rsidP=\"009168BC\" w:rsidR=\"00E02A2B\" w:rsidRDefault=\"00E02A2B\">' +\n    '<w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Test link before bookmark : </w:t></w:r>' +\n    fieldRunsXml +\n    '</w:p>';\n  pBefore.insertOoxml(wrapPackage(beforeXml), Word.InsertLocation.replace);\n}\n\nif (pBookmark) {\n  const bookmarkXml =\n    '<w:p w:rsidP=\"00C31A62\" w:rsidR=\"00C31A62\" w:rsidRDefault=\"00C31A62\">' +\n    '<w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Test bookmark : </w:t></w:r>' +\n    '<w:bookmarkStart w:name=\"bookmark1\" w:id=\"' + NEW_BMK_ID + '\"/>' +\n    '<w:r><w:t>bookmarked content</w:t></w:r>' +\n    '<w:bookmarkEnd w:id=\"' + NEW_BMK_ID + '\"/>' +\n    '</w:p>';\n  pBookmark.insertOoxml(wrapPackage(bookmarkXml), Word.InsertLocation.replace);\n}\n\nif (pAfter) {\n  const afterXml =\n    '<w:p w:rsidP=\"00E02A2B\" w:rsidR=\"00E02A2B\" w:rsidRDefault=\"00E02A2B\">' +\n    '<w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">Test link after bookmark : </w:t></w:r>' +\n    fieldRunsXml +\n    '<w:r w:rsidR=\"00D0546C\"><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '</w:p>';\n  pAfter.insertOoxml(wrapPackage(afterXml), Word.InsertLocation.replace);\n}\n\nawait context.sync();\nreturn \"ok\";\n", "ps1": "# The canonical XML diff for this change only touches two kinds of purely\n# internal/random identifiers that got regenerated when the fixture was\n# rebuilt going from 3.1.1 to 3.2.0:\n#   1. the (non-standard, tool-generated) w:rsidR=\"....\" GUID stamped on\n#      every run belonging to one of the two \"REF bookmark1\" field\n#      instances (begin/instrText/separate/result/end runs).\n#   2. the numeric w:id on the bookmarkStart/bookmarkEnd pair that wraps\n#      \"bookmarked content\" (bookmark \"bookmark1\").\n# No visible text, formatting or structure changes. We reproduce this by\n# rewriting, in place, the raw OOXML of the three paragraphs that contain\n# those runs/bookmark, via Range.InsertXML, using the exact literal target\n# values from the diff.\n\n$OLD_RSID = \"33AF09C2681E4769971B94CEE0D9A9EC\"\n$NEW_RSID = \"07822CC306A04453A57D7E97965A69F1\"\n$OLD_BMK_ID = \"45476227328674507871903955380577777817\"\n$NEW_BMK_ID = \"98058591413380703969263738417409166366\"\n\n$d = $word.ActiveDocument\n\nfunction Wrap-DocumentPackage($paragraphXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $paragraphXml + '</w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$fieldRunsXml = '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n    '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:instrText xml:space=\"preserve\"> REF bookmark1 \\h </w:instrText></w:r>' +\n    '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n    '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r>' +\n    '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r>'\n\n# Locate the three target paragraphs purely from their (unchanged) visible\n# text, so this does not depend on a fixed paragraph index.\n$pBefore = $null\n$pBookmark = $null\n$pAfter = $null\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text\n    if ($null -eq $pBefore -and $t -like \"*Test link before bookmark*\") {\n        $pBefore = $p\n    } elseif ($null -eq $pBookmark -and $t -like \"*Test bookmark*\" -and $t -like \"*bookmarked content*\" -and $t -notlike \"*duplicated*\") {\n        $pBookmark = $p\n    } elseif ($null -eq $pAfter -and $t -like \"*Test link after bookmark*\") {\n        $pAfter = $p\n    }\n}\n\nif ($pBefore) {\n    $beforeXml = '<w:p w:rsidP=\"009168BC\" w:rsidR=\"00E02A2B\" w:rsidRDefault=\"00E02A2B\">' +\n        '<w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr>' +\n        '<w:r><w:t xml:space=\"preserve\">Test link before bookmark : </w:t></w:r>' +\n        $fieldRunsXml +\n        '</w:p>'\n    [void]$pBefore.Range.InsertXML((Wrap-DocumentPackage $beforeXml))\n}\n\nif ($pBookmark) {\n    $bookmarkXml = '<w:p w:rsidP=\"00C31A62\" w:rsidR=\"00C31A62\" w:rsidRDefault=\"00C31A62\">' +\n        '<w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr>' +\n        '<w:r><w:t xml:space=\"preserve\">Test bookmark : </w:t></w:r>' +\n        '<w:bookmarkStart w:name=\"bookmark1\" w:id=\"' + $NEW_BMK_ID + '\"/>' +\n        '<w:r><w:t>bookmarked content</w:t></w:r>' +\n        '<w:bookmarkEnd w:id=\"' + $NEW_BMK_ID + '\"/>' +\n        '</w:p>'\n    [void]$pBookmark.Range.InsertXML((Wrap-DocumentPackage $bookmarkXml))\n}\n\nif ($pAfter) {\n    $afterXml = '<w:p w:rsidP=\"00E02A2B\" w:rsidR=\"00E02A2B\" w:rsidRDefault=\"00E02A2B\">' +\n        '<w:pPr><w:tabs><w:tab w:pos=\"3119\" w:val=\"left\"/></w:tabs></w:pPr>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n        '<w:r><w:t xml:space=\"preserve\">Test link after bookmark : </w:t></w:r>' +\n        $fieldRunsXml +\n        '<w:r w:rsidR=\"00D0546C\"><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n        '</w:p>'\n    [void]$pAfter.Range.InsertXML((Wrap-DocumentPackage $afterXml))\n}\n\nWrite-Output \"done\"\n"}
